$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the existing keyword text in B2 (remove trailing " 8888")
$ws.Range("B2").Value = "java.net.ConnectException: Connection timed out: connect"

# Add a new rule row (row 3), entering values in the same order the shared
# strings were originally authored in (KEYWORDS, then RULE, then ACTIONS)
$ws.Range("B3").Value = "com.automationanywhere.token.exception.SecurityTokenMissingException: UM1117.access.token.not.found"
$ws.Range("A3").Value = "R2"
$ws.Range("C3").Value = "1.Ask Client to clear all cache from browser and try to login again."

# Apply wrap-text formatting to the new row's B:C cells, matching row 2
$ws.Range("B3:C3").WrapText = $true

# Set row height for the new row to match row 2
$ws.Rows.Item(3).RowHeight = $ws.Rows.Item(2).RowHeight

# Update the active cell selection to C4, as happens after entering data through C3
$ws.Range("C4").Select()
